$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Yang Li"
$ws.Range("B3").Value = "Yes"

$ws.Range("B3").Select()
